# Fill in the missing "name" and "surname" values for rows 3 and 4,
# completing the client data pattern already present in the sheet.
# New shared strings must be appended in this order so the resulting
# sharedStrings.xml indexes match: 17 = "Имя 3", 18 = "Фамилия 2".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Имя 3"
$ws.Range("C3").Value = "Фамилия 2"

# Update the sheet's active selection, as recorded in the saved view state.
$ws.Range("C2:C3").Select()
